$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest crypto data snapshot.
# D-column values are forced to text (matching the sheets existing inlineStr storage)
# by setting a text NumberFormat before assignment, then resetting the cell style back
# to "Normal" so no residual formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.786.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.988.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.61%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.97%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.998.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.80%  "

$ws.Range("E10").Value = "  -3.83%  "

$ws.Range("E11").Value = "  -7.38%  "

$ws.Range("E12").Value = "  -3.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.509.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.70%  "

$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.840.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.993.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.93%  "

$ws.Range("E18").Value = "  -5.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  -3.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.115.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.76%  "

$ws.Range("E27").Value = "  -2.49%  "

$ws.Range("E28").Value = "  -3.43%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.26%  "

$ws.Range("E37").Value = "  -4.91%  "

$ws.Range("E38").Value = "  -4.84%  "

$ws.Range("E39").Value = "  -5.65%  "

$ws.Range("E40").Value = "  -7.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.416.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.78%  "

$ws.Range("E43").Value = "  -5.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.671"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("E49").Value = "  -3.62%  "

$ws.Range("E50").Value = "  -2.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.82%  "
